$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.823.54"
$ws.Range("E2").Value = "'  +1.30%  "
$ws.Range("D3").Value = "'3.315.11"
$ws.Range("E3").Value = "'  +6.31%  "
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("D5").Value = "'602.63"
$ws.Range("E5").Value = "'  +1.45%  "
$ws.Range("D6").Value = "'142.89"
$ws.Range("E6").Value = "'  +5.02%  "
$ws.Range("E7").Value = "'  +0.04%  "
$ws.Range("D8").Value = "'3.313.13"
$ws.Range("E8").Value = "'  +6.45%  "
$ws.Range("E9").Value = "'  +1.50%  "
$ws.Range("E10").Value = "'  +3.04%  "
$ws.Range("D11").Value = "'5.58"
$ws.Range("E11").Value = "'  +6.04%  "
$ws.Range("D12").Value = "'0.473"
$ws.Range("E12").Value = "'  +4.15%  "
$ws.Range("E13").Value = "'  +1.27%  "
$ws.Range("D14").Value = "'34.82"
$ws.Range("E14").Value = "'  +2.13%  "
$ws.Range("D15").Value = "'3.861.76"
$ws.Range("E15").Value = "'  +6.49%  "
$ws.Range("E16").Value = "'  +0.36%  "
$ws.Range("D17").Value = "'3.316.53"
$ws.Range("E17").Value = "'  +6.34%  "
$ws.Range("D18").Value = "'63.917.21"
$ws.Range("E18").Value = "'  +1.53%  "
$ws.Range("D19").Value = "'6.92"
$ws.Range("E19").Value = "'  +3.79%  "
$ws.Range("D20").Value = "'482.12"
$ws.Range("E20").Value = "'  +2.06%  "
$ws.Range("D21").Value = "'14.23"
$ws.Range("E21").Value = "'  +1.08%  "
$ws.Range("D22").Value = "'0.736"
$ws.Range("E22").Value = "'  +5.79%  "
$ws.Range("D23").Value = "'8.04"
$ws.Range("E23").Value = "'  +4.81%  "
$ws.Range("B24").Value = "'Litecoin"
$ws.Range("C24").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'84.98"
$ws.Range("E24").Value = "'  -1.16%  "
$ws.Range("B25").Value = "'InternetComputer(DFINITY)"
$ws.Range("C25").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "'13.47"
$ws.Range("E25").Value = "'  +4.35%  "
$ws.Range("E26").Value = "'  +0.07%  "
$ws.Range("E27").Value = "'  +1.95%  "
$ws.Range("D28").Value = "'7.27"
$ws.Range("E28").Value = "'  +4.93%  "
$ws.Range("E29").Value = "'  -0.10%  "
$ws.Range("E31").Value = "'  +4.29%  "
$ws.Range("D32").Value = "'29.38"
$ws.Range("E32").Value = "'  +10.20%  "
$ws.Range("E33").Value = "'  -1.32%  "
$ws.Range("E34").Value = "'  +1.18%  "
$ws.Range("E35").Value = "'  +2.72%  "
$ws.Range("D36").Value = "'5.98"
$ws.Range("E36").Value = "'  +3.35%  "
$ws.Range("D37").Value = "'53.01"
$ws.Range("E37").Value = "'  +1.94%  "
$ws.Range("D38").Value = "'0.0₃0751"
$ws.Range("E38").Value = "'  +7.81%  "
$ws.Range("E39").Value = "'  +4.63%  "
$ws.Range("D40").Value = "'431.47"
$ws.Range("E40").Value = "'  +3.44%  "
$ws.Range("D41").Value = "'3.048.53"
$ws.Range("E41").Value = "'  +5.12%  "
$ws.Range("D42").Value = "'8.41"
$ws.Range("E42").Value = "'  +2.85%  "
$ws.Range("D43").Value = "'2.75"
$ws.Range("E43").Value = "'  +2.50%  "
$ws.Range("E44").Value = "'  +0.10%  "
$ws.Range("E45").Value = "'  +0.18%  "
$ws.Range("E46").Value = "'  +4.61%  "
$ws.Range("D47").Value = "'26.44"
$ws.Range("E47").Value = "'  +3.86%  "
$ws.Range("D48").Value = "'35.97"
$ws.Range("E48").Value = "'  +14.85%  "
$ws.Range("E50").Value = "'  +2.25%  "
$ws.Range("D51").Value = "'2.32"
$ws.Range("E51").Value = "'  +3.17%  "
